$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Centre_Registration")
$custSearch = $wb.Worksheets.Item("Customer_CustSearch")

# --- New header cells E1:K1 (copy the existing header formatting from D1) ---
$ws.Range("D1").Copy()
$ws.Range("E1:K1").PasteSpecial(-4122)
$ws.Range("L1").PasteSpecial(-4122)
$ws.Range("M1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New body cells E2:K2 ---
# Text cells: copy the body formatting used elsewhere in the workbook (style
# matches Customer_CustSearch!J2) onto the text-valued cells.
$custSearch.Range("J2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("I2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Values are entered in the same order the author typed them so that the
# shared-string table indices line up with the source workbook.
$ws.Range("E1").Value = "StreetName"
$ws.Range("E2").Value = "2nd cross ashok nagar"
$ws.Range("F1").Value = "CentrName"
$ws.Range("F2").Value = "KHUSHI"
$ws.Range("G1").Value = "DisfromBranch"
$ws.Range("H1").Value = "BlockFld"
$ws.Range("H2").Value = "3rd"
$ws.Range("I1").Value = "LandMark"
$ws.Range("J1").Value = "Meeting_Time"
$ws.Range("I2").Value = "next to goeri appartments"
$ws.Range("K1").Value = "Meeting_Min"

# Numeric cells keep the default (unstyled) formatting.
$ws.Range("G2").Value = 10
$ws.Range("J2").Value = 10
$ws.Range("K2").Value = 30

# --- Column I width ---
$ws.Columns.Item(9).ColumnWidth = 13.5

# --- Selection moves to Q9 ---
[void]$ws.Range("Q9").Select()
